$d = $word.ActiveDocument

$replacements = @(
    @("2024-01-15 Monday", "2024-01-16 Tuesday"),
    @("681÷4=170, 1", "151÷5=30, 1"),
    @("681÷7=97, 2", "684÷6=114, 0"),
    @("463÷2=231, 1", "140÷4=35, 0"),
    @("504÷2=252, 0", "119÷6=19, 5"),
    @("525÷9=58, 3", "741÷6=123, 3"),
    @("404÷7=57, 5", "789÷5=157, 4"),
    @("649÷7=92, 5", "979÷7=139, 6"),
    @("408÷3=136, 0", "152÷5=30, 2"),
    @("916÷5=183, 1", "571÷4=142, 3"),
    @("785÷4=196, 1", "708÷4=177, 0"),
    @("689÷2=344, 1", "940÷4=235, 0"),
    @("577÷4=144, 1", "664÷5=132, 4"),
    @("403÷3=134, 1", "300÷9=33, 3"),
    @("850÷7=121, 3", "757÷8=94, 5"),
    @("937÷6=156, 1", "870÷5=174, 0"),
    @("876÷8=109, 4", "301÷3=100, 1"),
    @("988÷5=197, 3", "796÷5=159, 1"),
    @("620÷9=68, 8", "711÷3=237, 0"),
    @("565÷4=141, 1", "373÷4=93, 1"),
    @("700÷8=87, 4", "402÷4=100, 2"),
    @("357÷7=51, 0", "450÷5=90, 0"),
    @("864÷5=172, 4", "218÷5=43, 3"),
    @("235÷5=47, 0", "966÷3=322, 0"),
    @("490÷3=163, 1", "527÷5=105, 2"),
    @("744÷9=82, 6", "772÷5=154, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
